$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.808.83'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '2.290.70'
$ws.Range("E3").Value = '  -0.91%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.69'
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.63'
$ws.Range("E6").Value = '  -2.44%  '
$ws.Range("E7").Value = '  +0.57%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -3.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.61'
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("E11").Value = '  -0.38%  '
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.73'
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("E14").Value = '  -1.98%  '
$ws.Range("D15").Value = '2.647.19'
$ws.Range("E15").Value = '  -0.98%  '
$ws.Range("D16").Value = '2.292.99'
$ws.Range("E16").Value = '  +1.27%  '
$ws.Range("E17").Value = '  -1.86%  '
$ws.Range("D18").Value = '42.739.75'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("E19").Value = '  -4.54%  '
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("E21").Value = '  -2.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.77'
$ws.Range("E22").Value = '  -0.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.11'
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("E24").Value = '  -1.44%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("E26").Value = '  -1.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.02'
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.13'
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '165.88'
$ws.Range("E29").Value = '  -2.00%  '
$ws.Range("E30").Value = '  -0.81%  '
$ws.Range("E31").Value = '  -1.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.84'
$ws.Range("E32").Value = '  -1.67%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.83'
$ws.Range("E34").Value = '  -2.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.01'
$ws.Range("E35").Value = '  -3.48%  '
$ws.Range("E36").Value = '  -6.92%  '
$ws.Range("E37").Value = '  -1.23%  '
$ws.Range("E39").Value = '  -1.84%  '
$ws.Range("E40").Value = '  -3.63%  '
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.73'
$ws.Range("E42").Value = '  -0.78%  '
$ws.Range("D43").Value = '2.007.39'
$ws.Range("E43").Value = '  +0.71%  '
$ws.Range("E44").Value = '  -2.80%  '
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.07'
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.14'
$ws.Range("E47").Value = '  -1.64%  '
$ws.Range("E48").Value = '  -2.22%  '
$ws.Range("D49").Value = '2.514.49'
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.05'
$ws.Range("E50").Value = '  -3.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.97'
$ws.Range("E51").Value = '  -5.43%  '
